$wb = $excel.ActiveWorkbook

# --- Add the two new sheets (Areas, Buildings) ---
# Worksheets.Add() inserts the new sheet right before the active sheet, so
# adding Areas then Buildings (in that order) yields: Buildings, Areas, <rest>.
$areas = $wb.Worksheets.Add()
$areas.Name = "Areas"

$buildings = $wb.Worksheets.Add()
$buildings.Name = "Buildings"

# Reorder to the target order: Areas, Buildings, Carriers, ConsumerProducer, ...
# NOTE: worksheet handles returned by Add()/Item() are position-bound in this
# host, so always re-resolve by name immediately before each use rather than
# reusing a variable captured earlier in the script (it can silently resolve
# to a different sheet once the tab order changes).
$wb.Worksheets.Item("Buildings").Move($wb.Worksheets.Item("Carriers"))

# --- Populate Areas sheet headers (ID, Name first; then Scope, Parent_Area_ID,
# TopLevelArea, Area_WKT added - in this particular order - as the sheet was built) ---
$ws = $wb.Worksheets.Item("Areas")
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("D1").Value = "Scope"
$ws.Range("E1").Value = "Parent_Area_ID"
$ws.Range("C1").Value = "TopLevelArea"

# --- Populate Buildings sheet headers ---
$ws = $wb.Worksheets.Item("Buildings")
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "ESDLType"
$ws.Range("D1").Value = "Lat"
$ws.Range("E1").Value = "Lon"
$ws.Range("F1").Value = "Parent_Area_ID"

# --- Insert AreaBld_ID column into ConsumerProducer (new col F, before old Power col) ---
$ws = $wb.Worksheets.Item("ConsumerProducer")
$ws.Columns("F:F").Insert()
$ws = $wb.Worksheets.Item("ConsumerProducer")
$ws.Range("F1").Value = "AreaBld_ID"
$ws.Range("F2").Value = "NULL"
$ws.Range("F3").Value = "NULL"
$ws.Range("F4").Value = "NULL"

# --- Insert AreaBld_ID column into Conversion (new col F) ---
$ws = $wb.Worksheets.Item("Conversion")
$ws.Columns("F:F").Insert()
$ws = $wb.Worksheets.Item("Conversion")
$ws.Range("F1").Value = "AreaBld_ID"
$ws.Range("F2").Value = "NULL"
$ws.Range("F3").Value = "NULL"

# --- Insert AreaBld_ID column into Transport (new col F) ---
$ws = $wb.Worksheets.Item("Transport")
$ws.Columns("F:F").Insert()
$ws = $wb.Worksheets.Item("Transport")
$ws.Range("F1").Value = "AreaBld_ID"
$ws.Range("F2").Value = "NULL"

# --- Insert AreaBld_ID column into CablesPipesConnections (new col D) ---
$ws = $wb.Worksheets.Item("CablesPipesConnections")
$ws.Columns("D:D").Insert()
$ws = $wb.Worksheets.Item("CablesPipesConnections")
$ws.Range("D1").Value = "AreaBld_ID"
$ws.Range("D2").Value = "NULL"
$ws.Range("D3").Value = "NULL"
$ws.Range("D4").Value = "NULL"
$ws.Range("D5").Value = "NULL"
$ws.Range("D6").Value = "NULL"

# --- Finally, add the last Areas header (Area_WKT) ---
$wb.Worksheets.Item("Areas").Range("F1").Value = "Area_WKT"

# --- Make Areas the selected/active sheet (matches tabSelected moving off Carriers) ---
$wb.Worksheets.Item("Areas").Select()
